$d = $word.ActiveDocument

# 1) "Needs to update the other charts" -> strike-through (checked off)
$p1 = $d.Paragraphs(3)
$p1.Range.Font.StrikeThrough = 1

# 2) "Move the Cited vs Citing Chart bar ..." -> strike-through (checked off)
$p2 = $d.Paragraphs(7)
$p2.Range.Font.StrikeThrough = 1

# 3) Move the "_GoBack" bookmark from the end of the "Is there any way to make
#    this run more efficiently?" paragraph to the end of the "Make FDN svg ...
#    - not mine though" paragraph (right after the last run, before the
#    paragraph mark).
#
# The target position is immediately after the paragraph's last character,
# i.e. right up against the paragraph mark. Collapsed ranges placed exactly
# there confuse Bookmarks.Add, so we temporarily insert a placeholder
# character after that point, anchor the bookmark next to it, then remove
# the placeholder again; the bookmark stays put, collapsed, in the right spot.
$targetPara = $d.Paragraphs(10)
$endPos = $targetPara.Range.End - 1

$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$toRemove = $d.Range($endPos, $endPos + 1)
$toRemove.Text = ""
